$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 330.63635
$ws.Range("I33").Value = 141.28572
$ws.Range("J33").Value = 662
$ws.Range("K33").Value = 141.28572
$ws.Range("L33").Value = 662
$ws.Range("M33").Value = 87.71428
$ws.Range("N33").Value = -1120
$ws.Range("H97").Value = 1944.8334
$ws.Range("J97").Value = 1944.8334
$ws.Range("L97").Value = 5834.5002
$ws.Range("N97").Value = -6826.5002
$ws.Range("H112").Value = 2607.1482
$ws.Range("I112").Value = 995.8333
$ws.Range("K112").Value = 2987.4999
$ws.Range("M112").Value = -1879.4999
$ws.Range("H132").Value = 7581749
$ws.Range("I132").Value = 11910736
$ws.Range("J132").Value = 6020.875
$ws.Range("K132").Value = 35732208
$ws.Range("L132").Value = 18062.625
$ws.Range("M132").Value = -35729678
$ws.Range("N132").Value = -23122.625
$ws.Range("H135").Value = 893.8182
$ws.Range("I135").Value = 532.8889
$ws.Range("J135").Value = 2518
$ws.Range("K135").Value = 4796.0001
$ws.Range("L135").Value = 22662
$ws.Range("M135").Value = -2261.0001
$ws.Range("N135").Value = -27732
$ws.Range("H137").Value = 1011.86664
$ws.Range("I137").Value = 780.2381
$ws.Range("J137").Value = 1306.6666
$ws.Range("K137").Value = 2340.7143
$ws.Range("L137").Value = 3919.9998
$ws.Range("M137").Value = 209.2856999999999
$ws.Range("N137").Value = -9019.9998
$ws.Range("H138").Value = 1433.19
$ws.Range("I138").Value = 926.7778
$ws.Range("J138").Value = 1544.3536
$ws.Range("K138").Value = 2780.3334
$ws.Range("L138").Value = 4633.060799999999
$ws.Range("M138").Value = 2359.6666
$ws.Range("N138").Value = -14913.0608
$ws.Range("H141").Value = 597.7143
$ws.Range("I141").Value = 597.7143
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 1793.1429
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3386.8571
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1425.5135
$ws.Range("I61").Value = 1364.3214
$ws.Range("J61").Value = 1615.8889
$ws.Range("K61").Value = 1364.3214
$ws.Range("L61").Value = 1615.8889
$ws.Range("M61").Value = -1152.3214
$ws.Range("N61").Value = -2039.8889
$ws.Range("H74").Value = 2519
$ws.Range("I74").Value = 1041.3334
$ws.Range("J74").Value = 3405.6
$ws.Range("K74").Value = 1041.3334
$ws.Range("L74").Value = 3405.6
$ws.Range("M74").Value = -167.3334
$ws.Range("N74").Value = -5153.6
$ws.Range("H77").Value = 2519
$ws.Range("I77").Value = 1041.3334
$ws.Range("J77").Value = 3405.6
$ws.Range("K77").Value = 5206.666999999999
$ws.Range("L77").Value = 17028
$ws.Range("M77").Value = -838.6669999999995
$ws.Range("N77").Value = -25764
$ws.Range("H88").Value = 2436
$ws.Range("I88").Value = 1300
$ws.Range("J88").Value = 2720
$ws.Range("K88").Value = 1300
$ws.Range("L88").Value = 2720
$ws.Range("M88").Value = -894
$ws.Range("N88").Value = -3532
$ws.Range("H91").Value = 2436
$ws.Range("I91").Value = 1300
$ws.Range("J91").Value = 2720
$ws.Range("K91").Value = 1300
$ws.Range("L91").Value = 2720
$ws.Range("M91").Value = 104
$ws.Range("N91").Value = -5528
$ws.Range("H132").Value = 1803.4615
$ws.Range("I132").Value = 1536.8518
$ws.Range("K132").Value = 4610.555399999999
$ws.Range("M132").Value = -2080.555399999999
$ws.Range("H136").Value = 1425.5135
$ws.Range("I136").Value = 1364.3214
$ws.Range("J136").Value = 1615.8889
$ws.Range("K136").Value = 4092.9642
$ws.Range("L136").Value = 4847.6667
$ws.Range("M136").Value = -1542.9642
$ws.Range("N136").Value = -9947.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2604.9143
$ws.Range("I86").Value = 2902.5
$ws.Range("J86").Value = 1745.2222
$ws.Range("K86").Value = 2902.5
$ws.Range("L86").Value = 1745.2222
$ws.Range("M86").Value = -1779.5
$ws.Range("N86").Value = -3991.2222
$ws.Range("H89").Value = 2604.9143
$ws.Range("I89").Value = 2902.5
$ws.Range("J89").Value = 1745.2222
$ws.Range("K89").Value = 14512.5
$ws.Range("L89").Value = 8726.110999999999
$ws.Range("M89").Value = -8896.5
$ws.Range("N89").Value = -19958.111
$ws.Range("H134").Value = 3617.9363
$ws.Range("I134").Value = 847.79486
$ws.Range("J134").Value = 17122.375
$ws.Range("K134").Value = 2543.38458
$ws.Range("L134").Value = 51367.125
$ws.Range("M134").Value = -8.384579999999914
$ws.Range("N134").Value = -56437.125
$ws.Range("H138").Value = 53326.668
$ws.Range("J138").Value = 53326.668
$ws.Range("L138").Value = 53326.668
$ws.Range("N138").Value = -63606.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 90910550
$ws.Range("I16").Value = 125001256
$ws.Range("K16").Value = 125001256
$ws.Range("M16").Value = -125000969
$ws.Range("H31").Value = 1402.5111
$ws.Range("I31").Value = 1383.6666
$ws.Range("J31").Value = 1525
$ws.Range("K31").Value = 1383.6666
$ws.Range("L31").Value = 1525
$ws.Range("M31").Value = -1088.6666
$ws.Range("N31").Value = -2115
$ws.Range("H34").Value = 1402.5111
$ws.Range("I34").Value = 1383.6666
$ws.Range("J34").Value = 1525
$ws.Range("K34").Value = 1383.6666
$ws.Range("L34").Value = 1525
$ws.Range("M34").Value = -1181.6666
$ws.Range("N34").Value = -1929
$ws.Range("H58").Value = 676.59186
$ws.Range("I58").Value = 630.82355
$ws.Range("J58").Value = 780.3333
$ws.Range("K58").Value = 630.82355
$ws.Range("L58").Value = 780.3333
$ws.Range("M58").Value = -427.82355
$ws.Range("N58").Value = -1186.3333
$ws.Range("I94").Value = 528
$ws.Range("J94").Value = 599.2308
$ws.Range("K94").Value = 528
$ws.Range("L94").Value = 599.2308
$ws.Range("M94").Value = -77
$ws.Range("N94").Value = -1501.2308
$ws.Range("H105").Value = 1527.5
$ws.Range("I105").Value = 1536.6666
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 1536.6666
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 210.3334
$ws.Range("N105").Value = -4994
$ws.Range("H107").Value = 843.6
$ws.Range("I107").Value = 789.2
$ws.Range("J107").Value = 898
$ws.Range("K107").Value = 789.2
$ws.Range("L107").Value = 898
$ws.Range("M107").Value = 1130.8
$ws.Range("N107").Value = -4738
$ws.Range("H113").Value = 90910550
$ws.Range("I113").Value = 125001256
$ws.Range("K113").Value = 125001256
$ws.Range("M113").Value = -124999086
$ws.Range("H122").Value = 1783.5625
$ws.Range("I122").Value = 1237.1666
$ws.Range("K122").Value = 3711.4998
$ws.Range("M122").Value = -1261.4998
$ws.Range("H136").Value = 676.59186
$ws.Range("I136").Value = 630.82355
$ws.Range("J136").Value = 780.3333
$ws.Range("K136").Value = 1892.47065
$ws.Range("L136").Value = 2340.9999
$ws.Range("M136").Value = 657.5293500000002
$ws.Range("N136").Value = -7440.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 7107.3335
$ws.Range("J107").Value = 10462.4
$ws.Range("L107").Value = 31387.2
$ws.Range("N107").Value = -35227.2
$ws.Range("H121").Value = 840.875
$ws.Range("J121").Value = 999.5
$ws.Range("L121").Value = 2998.5
$ws.Range("N121").Value = -5618.5
$ws.Range("H131").Value = 34484170
$ws.Range("I131").Value = 76923310
$ws.Range("J131").Value = 2361.625
$ws.Range("K131").Value = 230769930
$ws.Range("L131").Value = 7084.875
$ws.Range("M131").Value = -230764890
$ws.Range("N131").Value = -17164.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2752.8823
$ws.Range("I80").Value = 1774.875
$ws.Range("J80").Value = 3622.2222
$ws.Range("K80").Value = 1774.875
$ws.Range("L80").Value = 3622.2222
$ws.Range("M80").Value = -776.875
$ws.Range("N80").Value = -5618.2222
$ws.Range("H83").Value = 2752.8823
$ws.Range("I83").Value = 1774.875
$ws.Range("J83").Value = 3622.2222
$ws.Range("K83").Value = 8874.375
$ws.Range("L83").Value = 18111.111
$ws.Range("M83").Value = -3882.375
$ws.Range("N83").Value = -28095.111
$ws.Range("H132").Value = 2354.182
$ws.Range("I132").Value = 1695.8572
$ws.Range("K132").Value = 5087.571599999999
$ws.Range("M132").Value = -2557.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2350
$ws.Range("I7").Value = 2400
$ws.Range("J7").Value = 2300
$ws.Range("K7").Value = 2400
$ws.Range("L7").Value = 2300
$ws.Range("M7").Value = -2288
$ws.Range("N7").Value = -2524
$ws.Range("H22").Value = 1000.2857
$ws.Range("I22").Value = 825
$ws.Range("J22").Value = 1234
$ws.Range("K22").Value = 825
$ws.Range("L22").Value = 1234
$ws.Range("M22").Value = -530
$ws.Range("N22").Value = -1824
$ws.Range("H27").Value = 1000.2857
$ws.Range("I27").Value = 825
$ws.Range("J27").Value = 1234
$ws.Range("K27").Value = 825
$ws.Range("L27").Value = 1234
$ws.Range("M27").Value = -718
$ws.Range("N27").Value = -1448
$ws.Range("H40").Value = 3227.0715
$ws.Range("I40").Value = 3046.75
$ws.Range("K40").Value = 3046.75
$ws.Range("M40").Value = -2910.75
$ws.Range("H46").Value = 5199.9
$ws.Range("J46").Value = 6857
$ws.Range("L46").Value = 6857
$ws.Range("N46").Value = -7233
$ws.Range("H122").Value = 50002800
$ws.Range("I122").Value = 50002800
$ws.Range("K122").Value = 150008400
$ws.Range("M122").Value = -150005950
$ws.Range("H126").Value = 2350
$ws.Range("I126").Value = 2400
$ws.Range("J126").Value = 2300
$ws.Range("K126").Value = 7200
$ws.Range("L126").Value = 6900
$ws.Range("M126").Value = -4730
$ws.Range("N126").Value = -11840
$ws.Range("H132").Value = 16020.464
$ws.Range("I132").Value = 1097.8182
$ws.Range("J132").Value = 42284.32
$ws.Range("K132").Value = 3293.4546
$ws.Range("L132").Value = 126852.96
$ws.Range("M132").Value = -763.4546
$ws.Range("N132").Value = -131912.96

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2399.7046
$ws.Range("I132").Value = 2434.5557
$ws.Range("J132").Value = 2242.875
$ws.Range("K132").Value = 7303.6671
$ws.Range("L132").Value = 6728.625
$ws.Range("M132").Value = -4773.6671
$ws.Range("N132").Value = -11788.625

